# excel2zugferd - TestRechnung.xlsx
# feat: #18 fixed and BT-73, BT-74, BT-134, BT-135 added
#
# Adds a new invoice line ("Anfahrt dazu" / travel expense, 0.5 h) to the
# "Rechnung2" sheet, which shifts the totals accordingly, and makes
# "Rechnung2" the active/selected sheet (instead of "Tabelle1").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Tabelle1
$ws2 = $wb.Worksheets.Item(2)   # Rechnung2

# --- New invoice position in row 20 --------------------------------------
# Pos. 11
$ws2.Range("A20").Value = 11

# Description - reuses the same look (bold-ish font, left aligned, wraps)
# as the other description cells in column C.
$descCell = $ws2.Range("C20")
$descCell.Value = "Anfahrt dazu"
$descCell.Font.Name = "Calibri"
$descCell.Font.Size = 11
$descCell.WrapText = $true
$descCell.HorizontalAlignment = -4131   # xlLeft

# Quantity: 0.5 (hours), shown with one decimal place
$qtyCell = $ws2.Range("D20")
$qtyCell.Value = 0.5
$qtyCell.NumberFormat = "#,##0.0"

# Unit: "h"
$ws2.Range("E20").Value = "h"

# Price / Sum formulas, consistent with the rest of the table
$ws2.Range("F20").Formula = '=IF(E20="10 Min.",22,75)'
$ws2.Range("G20").Formula = '=D20*F20'

# --- Make "Rechnung2" the active/selected sheet ---------------------------
$ws1.Range("A7").Select()
$ws2.Activate()

$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1
$ws2.Range("G36").Select()

$wb.Save()
